# Update generated output numbers (column F) on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet4 = $wb.Worksheets.Item("全部类型")

# --- "展览" sheet updates ---
$sheet1.Range("F2").Value  = 10083
$sheet1.Range("F4").Value  = 2521
$sheet1.Range("F11").Value = 1227
$sheet1.Range("F12").Value = 1042
$sheet1.Range("F13").Value = 3136
$sheet1.Range("F14").Value = 2348
$sheet1.Range("F16").Value = 2068
$sheet1.Range("F21").Value = 549
$sheet1.Range("F24").Value = 5
$sheet1.Range("F25").Value = 13
$sheet1.Range("F26").Value = 230
$sheet1.Range("F31").Value = 575
$sheet1.Range("F32").Value = 46
$sheet1.Range("F33").Value = 227
$sheet1.Range("F35").Value = 23
$sheet1.Range("F36").Value = 313
$sheet1.Range("F37").Value = 1654
$sheet1.Range("F38").Value = 105
$sheet1.Range("F41").Value = 435
$sheet1.Range("F42").Value = 937
$sheet1.Range("F44").Value = 345

# --- "全部类型" sheet updates ---
$sheet4.Range("F2").Value  = 10083
$sheet4.Range("F4").Value  = 2521
$sheet4.Range("F12").Value = 1227
$sheet4.Range("F13").Value = 1042
$sheet4.Range("F14").Value = 3136
$sheet4.Range("F15").Value = 2348
$sheet4.Range("F16").Value = 2068
$sheet4.Range("F17").Value = 2068
$sheet4.Range("F21").Value = 549
$sheet4.Range("F24").Value = 5
$sheet4.Range("F25").Value = 13
$sheet4.Range("F26").Value = 230
$sheet4.Range("F31").Value = 575
$sheet4.Range("F35").Value = 46
$sheet4.Range("F36").Value = 227
$sheet4.Range("F38").Value = 23
$sheet4.Range("F40").Value = 313
$sheet4.Range("F41").Value = 1654
$sheet4.Range("F42").Value = 105
$sheet4.Range("F46").Value = 435
$sheet4.Range("F47").Value = 937
$sheet4.Range("F49").Value = 345
